$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Data" - append three new rows of daily RRP data ---
$dataSheet = $wb.Worksheets.Item("Data")

$dataSheet.Range("A460").Value = 45142
$dataSheet.Range("B460").Value = 1793.804

$dataSheet.Range("A461").Value = 45145
$dataSheet.Range("B461").Value = 1810.583

$dataSheet.Range("A462").Value = 45146
$dataSheet.Range("B462").Value = 1778.351

# Match the existing date-cell styling (style index 3) used by column A
$dataSheet.Range("A460").Style = $dataSheet.Range("A459").Style
$dataSheet.Range("A461").Style = $dataSheet.Range("A459").Style
$dataSheet.Range("A462").Style = $dataSheet.Range("A459").Style

# --- Sheet 2: "SeriesInfo" - refresh the series metadata dates ---
$infoSheet = $wb.Worksheets.Item("SeriesInfo")

$infoSheet.Range("B3").Value = "2023-08-09"
$infoSheet.Range("B4").Value = "2023-08-09"
$infoSheet.Range("B7").Value = "2023-08-08"
$infoSheet.Range("B14").Value = "2023-08-08 13:01:06-05"
